$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value2 = 'мелочь'
$ws.Range("A16").Value2 = 'особливый товар'
$ws.Range("A17").Value2 = 'серебреный товар'
$ws.Range("A18").Value2 = 'деревенский товар'
$ws.Range("A19").Value2 = 'небогатый товар'
$ws.Range("A20").Value2 = 'крамными товар'
$ws.Range("A21").Value2 = 'мясо'
$ws.Range("A22").Value2 = 'железный товар'
$ws.Range("A24").Value2 = 'щепетильный товар'
$ws.Range("A25").Value2 = 'пушной товар'
$ws.Range("A26").Value2 = 'набойчатый товар'
$ws.Range("A27").Value2 = 'нужный товар'
$ws.Range("A29").Value2 = 'внутренний товар'
$ws.Range("A30").Value2 = 'питейный припасы'
$ws.Range("A31").Value2 = 'суровский товар'
$ws.Range("A32").Value2 = 'медный товар'
$ws.Range("A33").Value2 = 'привозный товар'
$ws.Range("A34").Value2 = 'оловянный товар'
$ws.Range("A37").Value2 = 'заморский товар'
$ws.Range("A38").Value2 = 'купецкий товар'
$ws.Range("A39").Value2 = 'домовый товар'
$ws.Range("A40").Value2 = 'харчевой припасы'
$ws.Range("A41").Value2 = 'меховой товар'
$ws.Range("A43").Value2 = 'надлежащий товар'
